$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.175.65'
$ws.Range("E2").Value = '  +0.36%  '

$ws.Range("D3").Value = '1.830.49'
$ws.Range("E3").Value = '  -0.23%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9994'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.20'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.18%  '

$ws.Range("E6").Value = '  +0.46%  '

$ws.Range("E7").Value = '  -0.13%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07348'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.35%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2908'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.56%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.22'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.78%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07633'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.70%  '

$ws.Range("D12").Value = '1.829.75'
$ws.Range("E12").Value = '  -0.22%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.974'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.43%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6703'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.22%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '82.37'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.28%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000008979'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.61%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.840'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.99%  '

$ws.Range("D18").Value = '29.160.57'
$ws.Range("E18").Value = '  +0.36%  '

$ws.Range("D19").Value = '2.078.58'
$ws.Range("E19").Value = '  -0.19%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '235.51'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.29%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.47'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.64%  '

$ws.Range("E22").Value = '  -0.19%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.345'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.40%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.003'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.17%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.52'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.43%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1388'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.43%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.519'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.48%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.61'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.05%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.490'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.45%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05853'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +6.38%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.228'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.79%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.083'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.50%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.085'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.61%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.859'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.47%  '

$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7255'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.30%  '

$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.138'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.11%  '

$ws.Range("E37").Value = '  -1.99%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.863'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.20%  '

$ws.Range("D39").Value = '1.229.81'
$ws.Range("E39").Value = '  +1.63%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01758'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.17%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.194'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.92%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9050'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.36%  '

$ws.Range("E43").Value = '  -0.01%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.90'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.04%  '

$ws.Range("D45").Value = '1.980.19'
$ws.Range("E45").Value = '  +0.10%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '65.82'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.67%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5042'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.99%  '

$ws.Range("E48").Value = '  -0.47%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.153'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.22%  '

$ws.Range("E50").Value = '  -4.18%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1131'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.90%  '
